# Rename the two existing sheets and add a third sheet, matching the
# "Added techlistic base, pages, tests" commit.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "LeftForm"

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "RightForm"

# Add the new sheet after the last existing sheet.
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "TechlisticForm"

# Seed the header row by copying the bold-header style from LeftForm so the
# new sheet reuses the existing "header" cell style instead of creating a
# brand-new one, then overwrite the text.
$ws1.Range("A1:E1").Copy($ws3.Range("A1:M1"))
$ws3.Rows.Item(1).RowHeight = 15.75

$headers = @("First Name ", "Last Name", "Gender", "Years of Exp", "Date", "Profession", "Automation Tools", "Continents", "Selenium Commands", "File Path", "Expected Results", "Actual Resuts", "Status")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws3.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$row2 = @("Arooba", "Imran", "Female", 1, "16-07-25", "Automation ", "Selenium ", "Asia", "WebElement Commands", "C:\Users\DELL\Desktop\dummy.txt", "All fields filled, Submit button clicked", "All fields filled, Submit button clicked", "PASS")
$row3 = @("abc", "xyz", "Male", 0, "16-07-25", "Manual ", "Selenium ", "Asia", "WebElement Commands", "C:\Users\DELL\Desktop\dummy.txt", "All fields filled, Submit button clicked", "All fields filled, Submit button clicked", "PASS")
$row4 = @("jkl", "mno", "Female", 1, "16-07-26", "Manual ", "Selenium ", "Asia", "WebElement Commands", "C:\Users\DELL\Desktop\dummy.txt", "All fields filled, Submit button clicked", "All fields filled, Submit button clicked", "PASS")

$rows = @($row2, $row3, $row4)
for ($r = 0; $r -lt $rows.Count; $r++) {
    $data = $rows[$r]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $ws3.Cells.Item($r + 2, $c + 1).Value = $data[$c]
    }
}

# Match column widths / "best fit" look of the other sheets.
$ws3.Columns.AutoFit()

# Selection the commit left behind on the new (now active) sheet.
$ws3.Range("J2:J4").Select()
$ws3.Range("J2").Activate()

# Selections left on the other two sheets by the edit.
$ws1.Range("C18").Select()
$ws2.Range("F5").Select()

# Re-apply best-fit widths on the renamed sheets too (content is unchanged,
# so this is a no-op size-wise but keeps behaviour consistent).
$ws1.Columns.AutoFit()
$ws2.Columns.AutoFit()

# TechlisticForm ends up the active/selected tab after the edit.
$ws3.Activate()
